$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Excel constant values used below ---
# xlContinuous = 1            (Borders.LineStyle - solid line)
# xlLineStyleNone = -4142     (Borders.LineStyle - no line)
# xlThin = 2                  (Borders.Weight)
# xlMedium = -4138            (Borders.Weight)
# xlHAlignCenter = -4108      (Range.HorizontalAlignment)
# xlVAlignBottom = -4107      (Range.VerticalAlignment - Excel default)
# xlPasteFormats = -4122      (Range.PasteSpecial)
# Borders.Item index: 7=Left, 8=Top, 9=Bottom, 10=Right

$xlContinuous = 1
$xlNone = -4142
$xlThin = 2
$xlMedium = -4138
$xlCenter = -4108
$xlBottom = -4107
$xlPasteFormats = -4122

# Clears all four edges of a cell's border (used before re-applying fresh
# borders on a cell whose format was copy/pasted from elsewhere).
function Clear-Borders($cell) {
    $cell.Borders.Item(7).LineStyle = $xlNone
    $cell.Borders.Item(8).LineStyle = $xlNone
    $cell.Borders.Item(9).LineStyle = $xlNone
    $cell.Borders.Item(10).LineStyle = $xlNone
}

# Applies a left/right/top/bottom border combo. Each side is "thin",
# "medium" or "none".
function Set-Border($cell, $left, $right, $top, $bottom) {
    if ($left -eq "thin") {
        $cell.Borders.Item(7).LineStyle = $xlContinuous
        $cell.Borders.Item(7).Weight = $xlThin
    } elseif ($left -eq "medium") {
        $cell.Borders.Item(7).LineStyle = $xlContinuous
        $cell.Borders.Item(7).Weight = $xlMedium
    }

    if ($right -eq "thin") {
        $cell.Borders.Item(10).LineStyle = $xlContinuous
        $cell.Borders.Item(10).Weight = $xlThin
    } elseif ($right -eq "medium") {
        $cell.Borders.Item(10).LineStyle = $xlContinuous
        $cell.Borders.Item(10).Weight = $xlMedium
    }

    if ($top -eq "thin") {
        $cell.Borders.Item(8).LineStyle = $xlContinuous
        $cell.Borders.Item(8).Weight = $xlThin
    } elseif ($top -eq "medium") {
        $cell.Borders.Item(8).LineStyle = $xlContinuous
        $cell.Borders.Item(8).Weight = $xlMedium
    }

    if ($bottom -eq "thin") {
        $cell.Borders.Item(9).LineStyle = $xlContinuous
        $cell.Borders.Item(9).Weight = $xlThin
    } elseif ($bottom -eq "medium") {
        $cell.Borders.Item(9).LineStyle = $xlContinuous
        $cell.Borders.Item(9).Weight = $xlMedium
    }
}

# Centers text horizontally while keeping Excel's default (bottom)
# vertical alignment, matching the target style's <alignment horizontal="center"/>.
function Set-CenterAlignment($cell) {
    $cell.VerticalAlignment = $xlBottom
    $cell.HorizontalAlignment = $xlCenter
}

# --- New small "Flow Rate" summary table (M6:O9) ---
# Cell values are assigned in the same order the original authoring
# session created them (Venturi, Pitot, then the merged header), so that
# new shared-string table entries land at the same indices as the target
# file (16=Venturi, 17=Pitot, 18=Flow Rate (m^3/s).

$m6 = $ws.Range("M6")
$n6 = $ws.Range("N6")
$o6 = $ws.Range("O6")
$m7 = $ws.Range("M7")
$n7 = $ws.Range("N7")
$o7 = $ws.Range("O7")
$m8 = $ws.Range("M8")
$n8 = $ws.Range("N8")
$o8 = $ws.Range("O8")
$m9 = $ws.Range("M9")
$n9 = $ws.Range("N9")
$o9 = $ws.Range("O9")

$n7.Value = "Venturi"
$o7.Value = "Pitot"
$n6.Value = "Flow Rate (m^3/s"
$m8.Value = "45 Hz"
$m9.Value = "47 Hz"
$n8.Value = 0.047899999999999998
$o8.Value = 0.21940000000000001
$n9.Value = 0.0516
$o9.Value = 0.22639999999999999

# Cells that use the lighter grey fill (fillId matching the sheet's "fill2")
# get their fill by copying the format of an existing cell using that
# fill (E6), then the borders are reset and rebuilt explicitly.
$fill2Source = $ws.Range("E6")
# Cells that use the darker grey fill ("fill3") copy their format from B3.
$fill3Source = $ws.Range("B3")

# Row 6 -----------------------------------------------------------------
$m6.Font.Size = 10
Set-CenterAlignment $m6

$fill2Source.Copy()
$n6.PasteSpecial($xlPasteFormats)
Clear-Borders $n6
Set-Border $n6 "medium" "thin" "medium" "medium"
Set-CenterAlignment $n6

$fill2Source.Copy()
$o6.PasteSpecial($xlPasteFormats)
Clear-Borders $o6
Set-Border $o6 "thin" "medium" "medium" "medium"
Set-CenterAlignment $o6

$ws.Range("N6:O6").Merge()

# Row 7 -----------------------------------------------------------------
$m7.Font.Size = 10
Set-CenterAlignment $m7

$fill3Source.Copy()
$n7.PasteSpecial($xlPasteFormats)
Clear-Borders $n7
Set-Border $n7 "medium" "thin" "none" "none"
Set-CenterAlignment $n7

$fill3Source.Copy()
$o7.PasteSpecial($xlPasteFormats)
Clear-Borders $o7
Set-Border $o7 "thin" "medium" "none" "none"
Set-CenterAlignment $o7

# Row 8 -------------------------------------------------------------------
$m8.Font.Size = 10
Set-Border $m8 "medium" "none" "medium" "thin"
Set-CenterAlignment $m8

$n8.Font.Size = 10
Set-Border $n8 "medium" "thin" "medium" "thin"
Set-CenterAlignment $n8

$o8.Font.Size = 10
Set-Border $o8 "thin" "medium" "medium" "thin"
Set-CenterAlignment $o8

# Row 9 -------------------------------------------------------------------
$fill2Source.Copy()
$m9.PasteSpecial($xlPasteFormats)
Clear-Borders $m9
Set-Border $m9 "medium" "none" "thin" "medium"
Set-CenterAlignment $m9

$fill2Source.Copy()
$n9.PasteSpecial($xlPasteFormats)
Clear-Borders $n9
Set-Border $n9 "medium" "thin" "thin" "medium"
Set-CenterAlignment $n9

$fill2Source.Copy()
$o9.PasteSpecial($xlPasteFormats)
Clear-Borders $o9
Set-Border $o9 "thin" "medium" "thin" "medium"
Set-CenterAlignment $o9
